$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Change the Runmode column (C) from "Y" to "N" for rows 2-13 and 15-18
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"
$ws.Range("C8").Value = "N"
$ws.Range("C9").Value = "N"
$ws.Range("C10").Value = "N"
$ws.Range("C11").Value = "N"
$ws.Range("C12").Value = "N"
$ws.Range("C13").Value = "N"
$ws.Range("C15").Value = "N"
$ws.Range("C16").Value = "N"
$ws.Range("C17").Value = "N"
$ws.Range("C18").Value = "N"

# Update the Results column (D2) from "SKIP" to "PASS"
$ws.Range("D2").Value = "PASS"

# Update the active selection to A31
$ws.Range("A31").Select()
